# "Generate Report for Handback" - update the localization-status report:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview sheet + each per-language sheet's Status column)
#   - Latest Handback DateTime refreshed per language
#   - The stale "handback file is not the latest" Error Detail is cleared
#     now that the handback is in sync
#   - Column widths are re-fit to the new (longer) Status text and the
#     now-empty Error Detail column

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column widths below are the "fit to new content" widths, expressed as the
# ColumnWidth (character units) that reproduces the target sheet column
# width for the new Status text / emptied Error Detail column.
$statusColWidth = 29.166666666666668
$errorDetailColWidth = 12.833333333333334

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Columns.Item(5).ColumnWidth = $statusColWidth
$ov.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("K2").Value = "2016-08-16 20:45:29"
$zh.Range("P2").Value = ""
$zh.Columns.Item(3).ColumnWidth = $statusColWidth
$zh.Columns.Item(16).ColumnWidth = $errorDetailColWidth

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("K2").Value = "2016-08-16 20:45:36"
$de.Range("P2").Value = ""
$de.Columns.Item(3).ColumnWidth = $statusColWidth
$de.Columns.Item(16).ColumnWidth = $errorDetailColWidth
